$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-off serial date/time value on the existing last row (row 10)
$ws.Range("A10").Value = 45810.39392270833

# Append new row 11 with the new price data point
$ws.Range("A11").Value = 45811.3937618142
$ws.Range("B11").Value = "EVOWHEY PROTEIN"
$ws.Range("C11").Value = "2Kg"
$ws.Range("D11").Value = "37,90€"

# Match the date/time number formatting used by the rest of column A
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat
